# Applies the weekly cryptos data refresh (GitHub Actions update).
# Most rows keep the same coin, but a handful of adjacent row-pairs
# (28/29, 42/43, 50/51) swap their ranking order along with new data.

function Set-TextValue($worksheet, $addr, $val) {
    # Assigning a numeric-looking string via .Value can make Excel
    # auto-convert it to a real number (and normalize/round it).
    # Force the cell to Text format first so the original text
    # representation (e.g. "0.717", "68.405.22") is preserved exactly,
    # then restore the cell's original style so formatting is unchanged.
    $cell = $worksheet.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" '68.405.22'
$ws.Range("E2").Value = '  -1.28%  '
Set-TextValue $ws "D3" '3.820.95'
$ws.Range("E3").Value = '  +1.98%  '
Set-TextValue $ws "D5" '599.68'
$ws.Range("E5").Value = '  -0.42%  '
Set-TextValue $ws "D6" '162.86'
$ws.Range("E6").Value = '  -3.37%  '
Set-TextValue $ws "D7" '3.818.66'
$ws.Range("E7").Value = '  +1.97%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -2.43%  '
$ws.Range("E10").Value = '  -3.40%  '
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("E13").Value = '  -4.00%  '
$ws.Range("E14").Value = '  -2.40%  '
Set-TextValue $ws "D15" '4.462.79'
$ws.Range("E15").Value = '  +2.04%  '
Set-TextValue $ws "D16" '3.781.83'
$ws.Range("E16").Value = '  +1.05%  '
Set-TextValue $ws "D17" '68.581.15'
$ws.Range("E17").Value = '  -0.96%  '
Set-TextValue $ws "D18" '7.54'
$ws.Range("E18").Value = '  +1.55%  '
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("E20").Value = '  -2.19%  '
Set-TextValue $ws "D21" '11.15'
$ws.Range("E21").Value = '  -1.45%  '
Set-TextValue $ws "D22" '484.46'
$ws.Range("E22").Value = '  -1.70%  '
Set-TextValue $ws "D23" '0.717'
$ws.Range("E23").Value = '  -1.80%  '
$ws.Range("E24").Value = '  +6.32%  '
Set-TextValue $ws "D25" '83.96'
$ws.Range("E25").Value = '  -1.05%  '
$ws.Range("E26").Value = '  -2.85%  '
Set-TextValue $ws "D27" '12.05'
$ws.Range("E27").Value = '  -2.18%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws "D28" '9.98'
$ws.Range("E28").Value = '  -0.89%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws "D29" '0.997'
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("E30").Value = '  -1.31%  '
Set-TextValue $ws "D31" '7.81'
$ws.Range("E31").Value = '  -4.67%  '
Set-TextValue $ws "D32" '3.974.28'
$ws.Range("E32").Value = '  +2.07%  '
Set-TextValue $ws "D33" '2.36'
$ws.Range("E33").Value = '  -4.43%  '
Set-TextValue $ws "D34" '31.69'
$ws.Range("E34").Value = '  +0.15%  '
Set-TextValue $ws "D35" '3.766.80'
$ws.Range("E35").Value = '  +2.38%  '
$ws.Range("E36").Value = '  -1.68%  '
$ws.Range("E37").Value = '  +0.93%  '
$ws.Range("E38").Value = '  -1.03%  '
Set-TextValue $ws "D39" '5.85'
$ws.Range("E39").Value = '  -2.32%  '
Set-TextValue $ws "D41" '0.317'
$ws.Range("E41").Value = '  -3.06%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws "D42" '431.93'
$ws.Range("E42").Value = '  +1.90%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws "D43" '2.94'
$ws.Range("E43").Value = '  -4.26%  '
Set-TextValue $ws "D44" '48.46'
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("E47").Value = '  -0.97%  '
Set-TextValue $ws "D48" '2.837.74'
$ws.Range("E48").Value = '  +1.53%  '
Set-TextValue $ws "D49" '142.63'
$ws.Range("E49").Value = '  +1.03%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws "D50" '0.0356'
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws "D51" '25.83'
$ws.Range("E51").Value = '  +12.79%  '
